$d = $word.ActiveDocument

# --- "KEY ACHIEVEMENTS AND IMPACT" section -----------------------------
# The same (or near-same) bullet text also appears earlier in the
# "PROFESSIONAL EXPERIENCE" section, so every Find must be scoped to just
# this section -- otherwise Find (which searches from the top of the
# document) would edit the wrong occurrence.
#
# Helper: build a fresh Range bounded by the "KEY ACHIEVEMENTS AND IMPACT"
# heading and the following "TECHNICAL SKILLS" heading, found by their
# (unique) heading text rather than a hard-coded paragraph index, so it
# keeps working as the section shrinks.

function Get-SectionRange() {
    $startPos = -1
    $endPos = -1
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $ptext = $d.Paragraphs.Item($i).Range.Text
        if ($startPos -eq -1 -and $ptext -like "*KEY ACHIEVEMENTS AND IMPACT*") {
            $startPos = $d.Paragraphs.Item($i).Range.End
        } elseif ($startPos -ne -1 -and $ptext -like "*TECHNICAL SKILLS*") {
            $endPos = $d.Paragraphs.Item($i).Range.Start
            break
        }
    }
    return $d.Range($startPos, $endPos)
}

function Replace-ExactTextInSection($find, $replace) {
    $scoped = Get-SectionRange
    $scoped.Find.Execute(
        $find,    # FindText
        $true,    # MatchCase
        $false,   # MatchWholeWord
        $false,   # MatchWildcards
        $false,   # MatchSoundsLike
        $false,   # MatchAllWordForms
        $true,    # Forward
        1,        # Wrap (wdFindContinue)
        $false,   # Format
        $replace, # ReplaceWith
        2         # Replace (wdReplaceOne)
    ) | Out-Null
}

function Delete-ParagraphWithTextInSection($find) {
    $scoped = Get-SectionRange
    $found = $scoped.Find.Execute(
        $find,   # FindText
        $true,   # MatchCase
        $false,  # MatchWholeWord
        $false,  # MatchWildcards
        $false,  # MatchSoundsLike
        $false,  # MatchAllWordForms
        $true,   # Forward
        1,       # Wrap (wdFindContinue)
        $false,  # Format
        "",      # ReplaceWith
        0        # Replace (wdReplaceNone)
    )
    if ($found) {
        # $scoped now wraps just the matched text (Find collapses it).
        # Extend by one character past the end to pull in the paragraph
        # mark too, so the whole paragraph -- not just its text -- goes
        # away (Paragraphs.Item(1).Range.Text is unreliable on a range
        # carved out of a larger Range, so operate on $scoped directly).
        $delRange = $d.Range($scoped.Start, $scoped.End + 1)
        $delRange.Delete()
    }
}

# 1) Rewrite four of the six bullets in place as impact statements.
Replace-ExactTextInSection `
    "• Built real-time FEC analysis systems using Python, Pandas and PySpark to detect likely fraud, money laundering and financial crimes across billions of records daily, performing time series analysis on trillions of records in the political spending sub-economy valued over `$2 trillion" `
    "• Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%"

Replace-ExactTextInSection `
    "• Built cloud-based data warehouse solutions on AWS processing billions of records with 99.94% accuracy" `
    "• `$4.7M savings enabled nonprofit access"

Replace-ExactTextInSection `
    "• Designed ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial datasets" `
    "• Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions"

Replace-ExactTextInSection `
    "• Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving classification accuracy from 23% to 64%" `
    "• 178% accuracy improvement in racial classification algorithms"

# 2) Remove the other two bullets outright (their content is no longer
#    represented in the section at all).
Delete-ParagraphWithTextInSection `
    "• Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M and enabling smaller nonprofits to conduct analysis"

Delete-ParagraphWithTextInSection `
    "• Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration, serving 12,847 analysts across 89 organizations"

Write-Output "Edit applied."
